# Update attendance/"want to go" counts (column F) across the
# "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# sheets, regenerated from the upstream data source.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 13
$wsExhibitions.Range("F4").Value = 13275
$wsExhibitions.Range("F9").Value = 114
$wsExhibitions.Range("F13").Value = 13238
$wsExhibitions.Range("F15").Value = 578
$wsExhibitions.Range("F16").Value = 8857
$wsExhibitions.Range("F17").Value = 7934
$wsExhibitions.Range("F18").Value = 233
$wsExhibitions.Range("F28").Value = 197
$wsExhibitions.Range("F29").Value = 120

# --- 演出 (Performances) sheet ---
$wsPerformances = $wb.Worksheets.Item("演出")
$wsPerformances.Range("F3").Value = 28

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 13
$wsAll.Range("F5").Value = 13275
$wsAll.Range("F10").Value = 114
$wsAll.Range("F14").Value = 13238
$wsAll.Range("F16").Value = 578
$wsAll.Range("F17").Value = 8857
$wsAll.Range("F18").Value = 7934
$wsAll.Range("F19").Value = 233
$wsAll.Range("F28").Value = 28
$wsAll.Range("F31").Value = 197
$wsAll.Range("F32").Value = 120
